# Generate Report for Handoff
# Updates the localization-status workbook to reflect that the
# a38ee76d-... resource has been handed off again ("Ready for handoff")
# with new "Latest Handoff Datetime" values on the per-locale sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "2016-03-04 05:59:19"

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("D3").Value = "2016-03-04 05:59:34"
